$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E2").Value = "2016-03-21 06:58:29"
$wsZh.Range("H2").Value = "2016-03-21 06:58:49"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E2").Value = "2016-03-21 06:58:32"
$wsDe.Range("H2").Value = "2016-03-21 06:58:54"
